$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from H1, the last populated header cell)
# onto the two new header cells I1 and J1, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for the new I / J columns, keyed by row number (2..37)
$values = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(6, 6)
    5  = @(8, 8)
    6  = @(6, 6)
    7  = @(13, 13)
    8  = @(8, 8)
    9  = @(6, 7)
    10 = @(7, 7)
    11 = @(6, 8)
    12 = @(6, 6)
    13 = @(2, 2)
    14 = @(9, 9)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(8, 9)
    18 = @(7, 7)
    19 = @(8, 9)
    20 = @(7, 8)
    21 = @(6, 6)
    22 = @(4, 5)
    23 = @(6, 7)
    24 = @(7, 8)
    25 = @(11, 11)
    26 = @(8, 8)
    27 = @(8, 8)
    28 = @(8, 8)
    29 = @(5, 5)
    30 = @(6, 6)
    31 = @(8, 9)
    32 = @(7, 7)
    33 = @(8, 8)
    34 = @(8, 8)
    35 = @(5, 5)
    36 = @(4, 4)
    37 = @(4, 4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}
